$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("C18").Select()

# Set Completion Date (column C) for rows 14-18 to 8/20/2020 (serial 44063).
# Copy the date number format from column B (style already used in the
# workbook) so no new cell style is created, then set the value.
for ($r = 14; $r -le 18; $r++) {
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("C$r").Value = 44063
}
$excel.CutCopyMode = $false

# Update Target Date (column B) for rows 19-35 (Excel serial date numbers;
# the cells already carry the date-formatted style, so this does not
# introduce any new cell style).
$ws.Range("B19").Value = 44063
$ws.Range("B20").Value = 44063
$ws.Range("B21").Value = 44063
$ws.Range("B22").Value = 44064
$ws.Range("B23").Value = 44064
$ws.Range("B24").Value = 44064
$ws.Range("B25").Value = 44064
$ws.Range("B26").Value = 44064
$ws.Range("B27").Value = 44065
$ws.Range("B28").Value = 44065
$ws.Range("B29").Value = 44065
$ws.Range("B30").Value = 44065
$ws.Range("B31").Value = 44066
$ws.Range("B32").Value = 44066
$ws.Range("B33").Value = 44067
$ws.Range("B34").Value = 44068
$ws.Range("B35").Value = 44069
